$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Quiz" column (column D). This shifts Vize/Fin/ORT one column left.
$ws.Range("D1").EntireColumn.Delete()

# Remove rows 6, 7 and 8 - only 4 data rows remain (rows 2-5).
$ws.Range("A6:A8").EntireRow.Delete()

# Update the values for the last remaining data row (row 5) with the new data.
$ws.Range("A5").Value = 210501014
$ws.Range("B5").Value = 52
$ws.Range("C5").Value = 63
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 77
$ws.Range("F5").Value = 99
